$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "Which Linux utility can be used to troubleshoot DNS problems?",
        "ques_type": 2,
        "options": [
            "dig",
            "ifconfig",
            "netstat",
            "telnet"
        ],
        "score": "dig"
    },
    {
        "title": "Which of the following are common continuous integration steps?",
        "ques_type": 15,
        "options": [
            "Compilation",
            "Deployment in a live system",
            "Code refactoring",
            "Unit testing",
            "Integration testing"
        ],
        "score": [
            "Compilation",
            "Unit testing",
            "Integration testing"
        ]
    },
    {
        "title": "Which of the following cloud concepts is best suited to enhance the availability of your application?",
        "ques_type": 2,
        "options": [
            "The use of cloud availability zones.",
            "The use of multiple cloud providers.",
            "The use of autoscaling groups.",
            "The use of hybrid cloud."
        ],
        "score": "The use of cloud availability zones."
    },
    {
        "title": "What git command will create a branch called feature/testgorilla?",
        "ques_type": 2,
        "options": [
            "git checkout -b feature/testgorilla",
            "git branch -c feature/testgorilla",
            "git create --branch feature/testgorilla",
            "git add --branch feature/testgorilla"
        ],
        "score": "git checkout -b feature/testgorilla"
    }
]
'@

# Remove the old row 2 (the cell holding the shared-string text) and strip
# the bold/bordered formatting that used to live on A1 (font/border no
# longer applied to any cell), then write the updated, pretty-printed text
# into the now-unstyled A1.
$ws.Range("A2").EntireRow.Delete() | Out-Null
$ws.Range("A1").ClearFormats() | Out-Null
$ws.Range("A1").Value = $newText
